$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha")
$ws.Range("E1:F1").EntireColumn.Insert()
$ws.Range("R1").EntireColumn.Insert()
$r2 = $ws.Range("R2")
$r2.ClearContents()
$r2.Borders.ColorIndex = 1
Write-Output "done"
